$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: FFR -> C ---
$ws.Range("C1").Value = "C"

# --- Row 2 (A Lag): update regression stats ---
$ws.Range("B2").Value = "-0.379***"

# C2 must become the literal text "0.485", not the number 0.485.
# Temporarily force a text number format so Excel stores it as a string,
# then clear the format back off so the cell doesn't keep a non-default
# style (matches the original workbook, where these cells carry no style).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.485"
$ws.Range("C2").ClearFormats()

# --- Row 3: label FFR Lag -> C Lag, update regression stats ---
$ws.Range("A3").Value = "C Lag"
$ws.Range("B3").Value = "-0.047***"
$ws.Range("C3").Value = "-0.785***"

# --- Remove old rows 4 (Constant) and 5 (r2_adj) entirely ---
$ws.Range("A4:C5").Delete()
